# New complement assay plot, new combined (DLS-Zeta-Pico) plot
# Rebuilds the "Sheet1" table: inserts an "N/P ratio" column, regroups rows
# by polymer name with the four N/P ratios (1, 5, 7.5, 10), and turns the
# former single "serum" row into the pDNA+H2O (N/P = 0) baseline row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# --- wipe the old A1:C23 table so we can lay the new A1:D22 table down clean
$ws.Range("A1:D23").ClearContents()

# --- header row ------------------------------------------------------------
$ws.Cells.Item(1,1).Value = "Polymer"
$ws.Cells.Item(1,2).Value = "N/P ratio"
$ws.Cells.Item(1,3).Value = "Zeta Potential"
$ws.Cells.Item(1,4).Value = "Stdev"

# --- data rows: Name, N/P ratio, Zeta Potential, Stdev ---------------------
$data = @(
    @("pDNA+" + [char]10 + "H" + [char]8322 + "O", 0,   -45,                4),
    @("S",   1,    -33.300000000000004,   0.26280537792569397),
    @("S",   5,    -24.4866666666667,     1.5545703229152699),
    @("S",   7.5,  -19.853333333333332,   2.5390330617951653),
    @("S",   10,   -21.78,                0.1512172829628507),
    @("B",   1,    -34.659999999999997,   0.59436240347675673),
    @("B",   5,    -32.54,                1.7439227811651148),
    @("B",   7.5,  -27.593333333333334,   0.89466318926298904),
    @("B",   10,   -27.816666666666666,   3.225381975656358),
    @("G1",  1,    -28.53,                2.7167750489627709),
    @("G1",  5,    -32.643333333333338,   1.569975229811676),
    @("G1",  7.5,  -26.149999999999995,   1.4236104336041757),
    @("G1",  10,   -24.056666666666668,   5.0125331808267335),
    @("G2",  1,    -28.150000000000002,   2.0909487479770199),
    @("G2",  5,    -25.416666666666668,   3.4069765025442775),
    @("G2",  7.5,  -28.386666666666667,   1.7841773703555608),
    @("G2",  10,   -18.846666666666668,   4.2736427345090728),
    @("G3",  1,    -31.353333333333335,   0.57702301128772637),
    @("G3",  5,    -25.180000000000003,   3.7113070473890915),
    @("G3",  7.5,  -20.970000000000002,   0.074833147735478819),
    @("G3",  10,   -15.780000000000001,   1.6678728968359673)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r,1).Value = $row[0]
    $ws.Cells.Item($r,2).Value = $row[1]
    $ws.Cells.Item($r,3).Value = $row[2]
    $ws.Cells.Item($r,4).Value = $row[3]
    $r = $r + 1
}

# --- formatting --------------------------------------------------------
# header "N/P ratio" -> Arial 10
$ws.Cells.Item(1,2).Font.Name = "Arial"
$ws.Cells.Item(1,2).Font.Size = 10

# pDNA+H2O label wraps onto two lines, row grows to fit
$ws.Cells.Item(2,1).WrapText = $true
$ws.Rows.Item(2).RowHeight = 32

# N/P ratio column: thin box border + Arial 10, first block centered horizontally
$top = $ws.Range("B2:B5")
$top.Font.Name = "Arial"
$top.Font.Size = 10
$top.Borders.LineStyle = 1
$top.HorizontalAlignment = -4108

# remaining N/P ratio cells: same border/font, but black font color and full
# center + wrap alignment
$rest = $ws.Range("B6:B22")
$rest.Font.Name = "Arial"
$rest.Font.Size = 10
$rest.Font.Color = 0
$rest.Borders.LineStyle = 1
$rest.HorizontalAlignment = -4108
$rest.VerticalAlignment = -4108
$rest.WrapText = $true

# --- sheet view bookkeeping ---------------------------------------------
$ws.Range("F18").Select()

# --- workbook window bookkeeping -----------------------------------------
$excel.ActiveWindow.Left = 11440
$excel.ActiveWindow.Top = 680
